# Update the "description" sheet (Sheet2) with a new validation note column
# and tweak the header text for the combranch_id note.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("description (รายละเอียด)")

# Update B1 header text (append extra note about numeric-only values)
$ws.Range("B1").Value = "เลขที่  (null)  ตัวเลขเท่านั้น"

# Widen column B slightly to fit the new text
$ws.Columns.Item(2).ColumnWidth = 25.1796875

# New columns H:I width
$ws.Columns.Item(8).ColumnWidth = 25.7265625
$ws.Columns.Item(9).ColumnWidth = 25.7265625

# Merge H1:I1 as an (empty) header cell matching style of other headers (border, centered)
$ws.Range("H1:I1").Merge()
$ws.Range("H1:I1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("H1:I1").VerticalAlignment = -4108    # xlCenter
$ws.Range("H1:I1").Borders.LineStyle = 1
$ws.Range("H1:I1").Borders.Weight = 2

# New content cells on row 2 describing validation notes
$ws.Range("H2").Value = "เป็นค่าว่างได้"
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("H2").VerticalAlignment = -4108

$ws.Range("I2").Value = "เพิ่มข้อมูล ต้องไม่เป็นค่าว่าง"
$ws.Range("I2").HorizontalAlignment = -4108

# Update selection to match recorded cursor position
$ws.Range("B4").Select()
